$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 516.6667
$ws.Range("I33").Value = 233.11111
$ws.Range("J33").Value = 1367.3334
$ws.Range("K33").Value = 233.11111
$ws.Range("L33").Value = 1367.3334
$ws.Range("M33").Value = -4.111109999999996
$ws.Range("N33").Value = -1825.3334
$ws.Range("H40").Value = 5842.857
$ws.Range("J40").Value = 7714.2856
$ws.Range("L40").Value = 7714.2856
$ws.Range("N40").Value = -8064.2856
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("H98").Value = 7700.143
$ws.Range("I98").Value = 7700.143
$ws.Range("K98").Value = 7700.143
$ws.Range("M98").Value = -6202.143
$ws.Range("H100").Value = 1349
$ws.Range("I100").Value = 1013
$ws.Range("K100").Value = 1013
$ws.Range("M100").Value = -472
$ws.Range("H112").Value = 2884.5715
$ws.Range("J112").Value = 2968
$ws.Range("L112").Value = 8904
$ws.Range("N112").Value = -11120
$ws.Range("H122").Value = 7700.143
$ws.Range("I122").Value = 7700.143
$ws.Range("K122").Value = 23100.429
$ws.Range("M122").Value = -20650.429
$ws.Range("H135").Value = 705.7143
$ws.Range("I135").Value = 529.2308
$ws.Range("K135").Value = 4763.077200000001
$ws.Range("M135").Value = -2228.077200000001
$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 2550.875
$ws.Range("I26").Value = 2584.5
$ws.Range("K26").Value = 2584.5
$ws.Range("M26").Value = -2254.5
$ws.Range("H61").Value = 7048.968
$ws.Range("I61").Value = 5948.28
$ws.Range("K61").Value = 5948.28
$ws.Range("M61").Value = -5736.28
$ws.Range("H74").Value = 3154.7837
$ws.Range("I74").Value = 2508.1072
$ws.Range("K74").Value = 2508.1072
$ws.Range("M74").Value = -1634.1072
$ws.Range("H77").Value = 3154.7837
$ws.Range("I77").Value = 2508.1072
$ws.Range("K77").Value = 12540.536
$ws.Range("M77").Value = -8172.536
$ws.Range("H136").Value = 7048.968
$ws.Range("I136").Value = 5948.28
$ws.Range("K136").Value = 17844.84
$ws.Range("M136").Value = -15294.84

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 11333
$ws.Range("J76").Value = 11333
$ws.Range("L76").Value = 11333
$ws.Range("N76").Value = -11963
$ws.Range("H79").Value = 11333
$ws.Range("J79").Value = 11333
$ws.Range("L79").Value = 11333
$ws.Range("N79").Value = -13517

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4819.852
$ws.Range("I31").Value = 3876.7058
$ws.Range("J31").Value = 6423.2
$ws.Range("K31").Value = 3876.7058
$ws.Range("L31").Value = 6423.2
$ws.Range("M31").Value = -3581.7058
$ws.Range("N31").Value = -7013.2
$ws.Range("H34").Value = 4819.852
$ws.Range("I34").Value = 3876.7058
$ws.Range("J34").Value = 6423.2
$ws.Range("K34").Value = 3876.7058
$ws.Range("L34").Value = 6423.2
$ws.Range("M34").Value = -3674.7058
$ws.Range("N34").Value = -6827.2
$ws.Range("H58").Value = 11161.429
$ws.Range("I58").Value = 7192.25
$ws.Range("J58").Value = 12749.1
$ws.Range("K58").Value = 7192.25
$ws.Range("L58").Value = 12749.1
$ws.Range("M58").Value = -6989.25
$ws.Range("N58").Value = -13155.1
$ws.Range("H95").Value = 17287.25
$ws.Range("J95").Value = 17287.25
$ws.Range("L95").Value = 17287.25
$ws.Range("N95").Value = -22779.25
$ws.Range("H134").Value = 8741.166999999999
$ws.Range("I134").Value = 7493.75
$ws.Range("K134").Value = 22481.25
$ws.Range("M134").Value = -19946.25
$ws.Range("H136").Value = 11161.429
$ws.Range("I136").Value = 7192.25
$ws.Range("J136").Value = 12749.1
$ws.Range("K136").Value = 21576.75
$ws.Range("L136").Value = 38247.3
$ws.Range("M136").Value = -19026.75
$ws.Range("N136").Value = -43347.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 12500
$ws.Range("J54").Value = 12500
$ws.Range("L54").Value = 37500
$ws.Range("N54").Value = -38618
$ws.Range("H131").Value = 23811728
$ws.Range("I131").Value = 62500850
$ws.Range("J131").Value = 3038.3845
$ws.Range("K131").Value = 187502550
$ws.Range("L131").Value = 9115.1535
$ws.Range("M131").Value = -187497510
$ws.Range("N131").Value = -19195.1535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 15107
$ws.Range("J15").Value = 15107
$ws.Range("L15").Value = 15107
$ws.Range("N15").Value = -15683
$ws.Range("H53").Value = 49999
$ws.Range("J53").Value = 49999
$ws.Range("L53").Value = 49999
$ws.Range("N53").Value = -51261
$ws.Range("H81").Value = 15107
$ws.Range("J81").Value = 15107
$ws.Range("L81").Value = 15107
$ws.Range("N81").Value = -17103
$ws.Range("H84").Value = 15107
$ws.Range("J84").Value = 15107
$ws.Range("L84").Value = 45321
$ws.Range("N84").Value = -55305
$ws.Range("H98").Value = 43333.332
$ws.Range("J98").Value = 43333.332
$ws.Range("L98").Value = 43333.332
$ws.Range("N98").Value = -49323.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3365
$ws.Range("I22").Value = 2581.25
$ws.Range("K22").Value = 2581.25
$ws.Range("M22").Value = -2286.25
$ws.Range("H27").Value = 3365
$ws.Range("I27").Value = 2581.25
$ws.Range("K27").Value = 2581.25
$ws.Range("M27").Value = -2474.25
$ws.Range("H70").Value = 25980
$ws.Range("J70").Value = 25980
$ws.Range("L70").Value = 25980
$ws.Range("N70").Value = -26520
$ws.Range("H73").Value = 25980
$ws.Range("J73").Value = 25980
$ws.Range("L73").Value = 25980
$ws.Range("N73").Value = -27852
$ws.Range("H100").Value = 2178791.5
$ws.Range("I100").Value = 3338136.5
$ws.Range("K100").Value = 3338136.5
$ws.Range("M100").Value = -3337595.5
$ws.Range("H116").Value = 69696
$ws.Range("J116").Value = 69696
$ws.Range("L116").Value = 69696
$ws.Range("N116").Value = -78874
$ws.Range("H132").Value = 7861.5415
$ws.Range("I132").Value = 8069.6665
$ws.Range("K132").Value = 24208.9995
$ws.Range("M132").Value = -21678.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 40000
$ws.Range("J58").Value = 40000
$ws.Range("L58").Value = 40000
$ws.Range("N58").Value = -40616
$ws.Range("H132").Value = 2834.8647
$ws.Range("I132").Value = 2682.5715
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 8047.7145
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -5517.7145
$ws.Range("N132").Value = -21560
$ws.Range("H136").Value = 5029.5
$ws.Range("I136").Value = 3857.6
$ws.Range("K136").Value = 11572.8
$ws.Range("M136").Value = -9022.799999999999
